# Adjustments to scenario 3a and 3b parameters for STH
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Platform Coverage" ---
$ws1 = $wb.Worksheets.Item("Platform Coverage")

# Extend row 2 (All/Treatment/Campaign/MDA, age 5-15) coverage of 0.6
# across the remaining biennial year columns (previously only set on row 3).
$ws1.Range("P2").Value = 0.6
$ws1.Range("R2").Value = 0.6
$ws1.Range("T2").Value = 0.6
$ws1.Range("V2").Value = 0.6
$ws1.Range("X2").Value = 0.6
$ws1.Range("Z2").Value = 0.6
$ws1.Range("AB2").Value = 0.6
$ws1.Range("AD2").Value = 0.6

# Remove row 3 (the duplicate All/Treatment/Campaign/MDA row with the 0.7
# coverage values) - all rows below shift up by one.
$ws1.Rows.Item(3).Delete()

# --- Sheet 2: "MarketShare" ---
$ws2 = $wb.Worksheets.Item("MarketShare")

# Clear the New Product A market share values (row 2, years 2026-2040)
$ws2.Range("L2:Z2").ClearContents()

# Old Product B (SOC) now holds 100% market share for the full timeline
$ws2.Range("L3:Z3").Value = 1

# Update the cell selections left on each sheet (select sheet 2's cell
# first so that sheet 1 / "Platform Coverage" remains the active tab,
# matching the original tabSelected state).
$ws2.Range("Z3").Select()
$ws1.Range("G6").Select()
